$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TTD")
$ws.Activate()

$ws.Range("A2").Value = "Add: Required`nEdit: Required`nEdit Rates: Required`nRetrieve Batch: Optional`nRetrieve Custom: Required"
$ws.Range("B2").Value = "Add: Required`nEdit: Required`nEdit Rates: Optional`nRetrieve Batch: Optional`nRetrieve Custom: Optional"
$ws.Range("C2").Value = "Add: Required`nEdit: Required`nEdit Rates: Optional`nRetrieve Batch: Optional`nRetrieve Custom: Optional"
$ws.Range("D2").Value = "Add: Required`nEdit: Required`nEdit Rates: Optional`nRetrieve Batch: Optional`nRetrieve Custom: Optional"
$ws.Range("E2").Value = "Add: Required`nEdit: Required`nEdit Rates: Optional`nRetrieve Batch: Optional`nRetrieve Custom: Optional"
$ws.Range("F2").Value = "Add: Required`nEdit: Not Required`nEdit Rates: Required`nRetrieve Batch: Optional`nRetrieve Custom: Optional`n`"bombora`" or `"eyeota`" only"
$ws.Range("G2").Value = "Add: Required`nEdit: Not Required`nEdit Rates: Required`nRetrieve Batch: Optional`nRetrieve Custom: Required"
$ws.Range("H2").Value = "Add: Required`nEdit: Not Required`nEdit Rates: Required`nRetrieve Batch: Optional`nRetrieve Custom: Optional"
$ws.Range("I2").Value = "Add: Required`nEdit: Not Required`nEdit Rates: Required`nRetrieve Batch: Optional`nRetrieve Custom: Optional`nValues: CPM or PercentOfMediaCost"
$ws.Range("J2").Value = "Add: Not Required`nEdit: Not Required`nEdit Rates: Optional`nRetrieve Batch: Required`nRetrieve Custom: Optional"

$ws.Columns.Item(1).ColumnWidth = 22.5
$ws.Columns.Item(7).ColumnWidth = 23.0
$ws.Columns.Item(8).ColumnWidth = 23.0
$ws.Columns.Item(9).ColumnWidth = 23.333333333333336
$ws.Columns.Item(10).ColumnWidth = 23.833333333333336

$ws.Range("D7").Select()
